# Actualización automática 2025-07-04 14:10:07
#
# Applies the updated "julio" sales figures for advisor CASTRO ALCIVAR EDA MARIA:
#   - PORCEKER S.A. now has 518.4 in "240X120 PORCELANATO"
#   - SALAZAR BALLADARES MARIA ANGELICA has a -21.42 adjustment in "PORCELANATO"
# and propagates the change through the monthly and compliance summary sheets.

$wb = $excel.ActiveWorkbook

$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# --- Sheet "VENTAS POR GRUPO" -------------------------------------------------
# C36: PORCEKER S.A. / 240X120 PORCELANATO
$wsGrupo.Range("C36").Value = 518.4
# M43: SALAZAR BALLADARES MARIA ANGELICA / PORCELANATO
$wsGrupo.Range("M43").Value = -21.42
# C55: running count of non-zero clients for the 240X120 PORCELANATO column
$wsGrupo.Range("C55").Value = "1 de 53"

# --- Sheet "VENTA MENSUAL" ----------------------------------------------------
# F36: PORCEKER S.A. / julio
$wsMensual.Range("F36").Value = 518.4
# F43: SALAZAR BALLADARES MARIA ANGELICA / julio
$wsMensual.Range("F43").Value = 4.68
# F55: julio column total
$wsMensual.Range("F55").Value = 6876.440000000001

# --- Sheet "CUMPLIMIENTO MENSUAL" --------------------------------------------
# Row 2: 240X120 PORCELANATO
$wsCumplimiento.Range("D2").Value = 518.4
$wsCumplimiento.Range("E2").Value = 5301.6
$wsCumplimiento.Range("F2").Value = 0.0890721649484536

# Row 16: PORCELANATO
$wsCumplimiento.Range("D16").Value = 2454.86
$wsCumplimiento.Range("E16").Value = 52266.37
$wsCumplimiento.Range("F16").Value = 0.0448611992091552

# Row 19: TOTAL
$wsCumplimiento.Range("D19").Value = 7284.27
$wsCumplimiento.Range("E19").Value = 97928.59999999999
$wsCumplimiento.Range("F19").Value = 0.06923364033316458
